$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - date advances to 45946; B, C, F cleared (no new forecast yet); D and E updated
$ws.Range("A2").Value = 45946
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 6644
$ws.Range("E2").Value = 8231.102329
$ws.Range("F2").ClearContents()

# Row 3
$ws.Range("A3").Value = 45947
$ws.Range("B3").Value = 6014.75591687487
$ws.Range("C3").Value = 5545.35648821073
$ws.Range("D3").Value = 3620
$ws.Range("E3").Value = 8901.479227
$ws.Range("F3").Value = 200.503324930661

# Row 4
$ws.Range("A4").Value = 45948
$ws.Range("B4").Value = 1986.17410438121
$ws.Range("C4").Value = 4067.96160906929
$ws.Range("D4").Value = 3620
$ws.Range("E4").Value = 4586.915774
$ws.Range("F4").Value = 127.02930327867

# Row 5
$ws.Range("A5").Value = 45949
$ws.Range("B5").Value = 1967.69496430682
$ws.Range("C5").Value = 4178.689826757
$ws.Range("D5").Value = 3620
$ws.Range("E5").Value = 4638.699052
$ws.Range("F5").Value = 134.570579768758

# Row 6
$ws.Range("A6").Value = 45950
$ws.Range("B6").Value = 7100.81174419054
$ws.Range("C6").Value = 6936.94748718724
$ws.Range("D6").Value = 3620
$ws.Range("E6").Value = 10748.9061
$ws.Range("F6").Value = 290.210076791529

# Row 7
$ws.Range("A7").Value = 45951
$ws.Range("B7").Value = 6392.95297294923
$ws.Range("C7").Value = 6485.73300313033
$ws.Range("D7").Value = 3620
$ws.Range("E7").Value = 9928.5279
$ws.Range("F7").Value = 266.721163757546

# Row 8
$ws.Range("A8").Value = 45952
$ws.Range("B8").Value = 6392.95297294923
$ws.Range("C8").Value = 6389.75480455025
$ws.Range("D8").Value = 3620
$ws.Range("E8").Value = 9928.5279
$ws.Range("F8").Value = 262.722072150043

# Row 9
$ws.Range("A9").Value = 45953
$ws.Range("B9").Value = 6392.95297294923
$ws.Range("C9").Value = 6119.72960869775
$ws.Range("D9").Value = 3620
$ws.Range("E9").Value = 9928.5279
$ws.Range("F9").Value = 251.471022322855

# Row 10
$ws.Range("A10").Value = 45954
$ws.Range("B10").Value = 6392.95297294923
$ws.Range("C10").Value = 5464.78991455457
$ws.Range("D10").Value = 3620
$ws.Range("E10").Value = 9928.5279
$ws.Range("F10").Value = 224.181868400222

# Row 11
$ws.Range("A11").Value = 45955
$ws.Range("B11").Value = 2044.14689306201
$ws.Range("C11").Value = 3632.02424946153
$ws.Range("D11").Value = 3620
$ws.Range("E11").Value = 5157.7877
$ws.Range("F11").Value = 130.236044016647

# Row 12
$ws.Range("A12").Value = 45956
$ws.Range("B12").Value = 1933.9611505005
$ws.Range("C12").Value = 3525.03286882277
$ws.Range("D12").Value = 3620
$ws.Range("E12").Value = 5037.340811
$ws.Range("F12").Value = 125.350522055094

# Row 13
$ws.Range("A13").Value = 45957
$ws.Range("B13").Value = 6757.64823424925
$ws.Range("C13").Value = 5927.46358646827
$ws.Range("D13").Value = 3620
$ws.Range("E13").Value = 10717.733504
$ws.Range("F13").Value = 261.147869009126

# Row 14
$ws.Range("A14").Value = 45958
$ws.Range("B14").Value = 6757.64823424925
$ws.Range("C14").Value = 6000.26642188892
$ws.Range("D14").Value = 3620
$ws.Range("E14").Value = 10717.733504
$ws.Range("F14").Value = 264.181320484986

# Row 15
$ws.Range("A15").Value = 45959
$ws.Range("B15").Value = 6757.64823424925
$ws.Range("C15").Value = 5939.40109343239
$ws.Range("D15").Value = 3620
$ws.Range("E15").Value = 10717.733504
$ws.Range("F15").Value = 261.645265132631
